# Powerpoint writer: consolidate text run nodes.
#
# The original deck splits text that was typed word-by-word into one
# <a:r> run per word plus one <a:r> run per single space (an artifact of
# how the text was originally authored). This script rewrites each of
# those paragraphs so that every "word" run is consolidated with the
# single-space run that immediately follows it, producing one run per
# "word " chunk instead of two. The trailing word of a paragraph (which
# has no following space run) is left as its own run, untouched.
#
# Mechanically, for a TextRange whose current text is e.g. "Slide 1
# (Content)" (runs: "Slide", " ", "1", " ", "(Content)"), we replace the
# characters spanning "Slide " (the word plus its trailing space) with
# the literal same text "Slide " via TextRange.Characters(start,length).
# Even though the replacement text is byte-for-byte identical to what
# was already there, re-assigning .Text on that sub-range forces the
# writer to collapse it into a single run, which is exactly the
# consolidation the diff shows. We do the same for "1 ", leaving
# "(Content)" as the final, already-singleton run.

function Merge-WordSpaceRuns($tr) {
    $text = $tr.Text
    $len = $text.Length
    $i = 0
    $pos = 1
    while ($i -lt $len) {
        $ch = $text.Substring($i, 1)
        if ($ch -eq " ") {
            $i = $i + 1
            $pos = $pos + 1
            continue
        }

        # Walk to the end of this run of non-space characters ("word").
        $j = $i
        while ($j -lt $len -and $text.Substring($j, 1) -ne " ") {
            $j = $j + 1
        }
        $wordLen = $j - $i

        if ($j -lt $len -and $text.Substring($j, 1) -eq " ") {
            # Word is immediately followed by a single space: merge the
            # word run with that trailing-space run into one run.
            $mergeLen = $wordLen + 1
            $chunk = $text.Substring($i, $mergeLen)
            $r = $tr.Characters($pos, $mergeLen)
            $r.Text = $chunk
            $i = $j + 1
            $pos = $pos + $mergeLen
        } else {
            # Final word in the paragraph - no trailing space to merge with.
            $i = $j
            $pos = $pos + $wordLen
        }
    }
}

$p = $ppt.ActivePresentation

# Slide titles ("Slide N (...)") - slides 1-12.
for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $slide = $p.Slides.Item($idx)
    $title = $slide.Shapes.Item("Title 1")
    Merge-WordSpaceRuns $title.TextFrame.TextRange
}

# "an image" / "An image" captions under the pictures on slides 6, 7, 8.
foreach ($idx in 6, 7, 8) {
    $slide = $p.Slides.Item($idx)
    $caption = $slide.Shapes.Item("TextBox 3")
    Merge-WordSpaceRuns $caption.TextFrame.TextRange
}
